# "Updated to remove anchor underlines"
#
# - Rename Sheet1 -> "Personal Access Tokens"
# - Remove the stored GitHub PAT row (A3/B3: token + "docs workflow"),
#   leaving only the "Token" / "Purpose" header row
# - Move the active-cell selection to A26
# - Set the page orientation (adds the pageSetup element on export)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Personal Access Tokens"

# Drop the secret-token data row entirely so the used range shrinks back
# to just the header row (A1:B1) and the shared strings for the token
# value + "docs workflow" purpose are dropped along with it.
$ws.Range("A3:B3").Delete()

# Match the saved selection/active cell.
$ws.Range("A26").Select() | Out-Null

# Touch PageSetup so a <pageSetup> element is emitted on save.
$ws.PageSetup.Orientation = 1
